# Applies the "Updated symbol list" data refresh to Sheet1 of cryptos.xlsx.
# Column D (Price) and G (Hora) values are stored as text in the workbook,
# so numeric-looking replacement values are written via a text number format
# and the cell style is then reset back to Normal (no lingering formatting).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Cells whose new content is numeric-looking text (Price/Hora columns).
$textNumericUpdates = @{
    "D2" = "275.03"
    "G2" = "4"
    "D3" = "23.15"
    "G3" = "4"
    "D4" = "6.476"
    "G4" = "4"
    "D5" = "0.06295"
    "G5" = "4"
    "D6" = "3.662"
    "G6" = "4"
    "D7" = "6.690"
    "G7" = "4"
    "D8" = "1.397"
    "G8" = "4"
    "D9" = "0.8347"
    "G9" = "4"
    "D10" = "0.01387"
    "G10" = "4"
    "D11" = "0.1628"
    "G11" = "4"
    "D12" = "0.08279"
    "G12" = "4"
    "D13" = "0.03432"
    "G13" = "4"
    "D14" = "0.03108"
    "G14" = "4"
    "D15" = "0.09299"
    "G15" = "4"
    "D16" = "3.842"
    "G16" = "4"
    "D17" = "0.001641"
    "G17" = "4"
    "D18" = "0.04789"
    "G18" = "4"
    "D19" = "0.006371"
    "G19" = "4"
    "D20" = "0.005687"
    "G20" = "4"
    "G21" = "4"
    "G22" = "4"
    "D23" = "3.713"
    "G23" = "4"
    "D24" = "2.409"
    "G24" = "4"
    "G25" = "4"
    "G26" = "4"
    "G27" = "4"
    "G28" = "4"
    "G29" = "4"
    "G30" = "4"
    "G31" = "4"
    "G32" = "4"
    "G33" = "4"
    "G34" = "4"
    "G35" = "4"
    "G36" = "4"
    "G37" = "4"
    "G38" = "4"
    "G39" = "4"
    "D40" = "0.04722"
    "G40" = "4"
    "D41" = "0.007056"
    "G41" = "4"
    "G42" = "4"
    "D43" = "0.003300"
    "G43" = "4"
    "D44" = "0.01214"
    "G44" = "4"
    "D45" = "0.00006245"
    "G45" = "4"
    "G46" = "4"
    "G47" = "4"
    "D48" = "0.7966"
    "G48" = "4"
    "D49" = "0.01594"
    "G49" = "4"
    "D50" = "0.00002300"
    "G50" = "4"
    "G51" = "4"
}
foreach ($cellRef in $textNumericUpdates.Keys) {
    Set-TextValue $cellRef $textNumericUpdates[$cellRef]
}

# Row 49/50 coin name/link/volume swap (non-numeric text, safe to set directly).
$plainTextUpdates = @{
    "B49" = "BOLO"
    "C49" = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
    "E49" = "48BOLOBOLOBestin24h"
    "B50" = "CryptobidCoin"
    "C50" = "https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc"
    "E50" = "49CryptobidCoinCBC"
}
foreach ($cellRef in $plainTextUpdates.Keys) {
    $ws.Range($cellRef).Value = $plainTextUpdates[$cellRef]
}
